$d = $word.ActiveDocument

# Locate the paragraph that currently contains only "Lab" (the last
# paragraph in the document, per the diff context).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r", "`n", "`v") -eq "Lab") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Lab' paragraph to extend."
}

$r = $target.Range
# Exclude the trailing paragraph mark so new runs are appended inside
# the same paragraph, right after the existing "Lab" run.
$r.End = $r.End - 1

$newContentXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t xml:space="preserve"> 28 &amp; 29 &#8211; Mitral </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>Stenosis</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> and Aortic Regurgitation -- .DES files for Aortic and Mitral </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>Stenosis</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> and Regurgitation under Display are changed in order to scale by .01 mm^2 rather than .1 mm^2.</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r.InsertXML($newContentXml)
